# Splits the combined "To update the Customer's..." / "If the reservation
# balance is now $0.00..." sentence pair in the Admin Screen TL;DR bullet into
# the new wording ("...enter any additional payment you wish to make on the
# reservation. If the balance is paid in full, enter 0..") while keeping the
# existing run formatting (Arial Nova Light) and adding the Word grammar-check
# markers (<w:proofErr>) the author's edit introduced around "0..".

$d = $word.ActiveDocument

$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "To update the Customer’s reservation balance, enter the current payment amount and select the ‘Save’ option. If the reservation balance is now $0.00, a new ticket is generated. ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target sentence to replace."
}

# Re-seat a fresh Range over the exact match bounds: calling InsertXML
# directly on the Find-mutated range object appends instead of replacing,
# but a freshly constructed Range over the same Start/End replaces in place.
$target = $d.Range($searchRange.Start, $searchRange.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t xml:space="preserve">To update the </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t xml:space="preserve">Customer’s reservation balance, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t>enter any additional payment you wish to make on the reservation</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t xml:space="preserve"> If the balance is paid in full, enter </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t>0.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
